$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: convert C2:G2 from numeric to text, keep values the same
$ws.Range("C2:G2").NumberFormat = "@"
$ws.Range("C2").Value = "14.8881"
$ws.Range("D2").Value = "120.7855"
$ws.Range("E2").Value = "4663"
$ws.Range("F2").Value = "92"
$ws.Range("G2").Value = "1000"

# Row 3: convert C3:G3 from numeric to text, G3 value changes from 2 to 500
$ws.Range("C3:G3").NumberFormat = "@"
$ws.Range("C3").Value = "14.9"
$ws.Range("D3").Value = "120.78"
$ws.Range("E3").Value = "3"
$ws.Range("F3").Value = "3"
$ws.Range("G3").Value = "500"
